# Update EMMOntoPy onto_update.xlsx template to use emmo 1.0.0-beta4
# instead of the beta/emmo-inferred-chemistry2 version, rename the
# example "Pattern" concept to "SpecialPattern", turn the imported
# ontology URL into a real hyperlink, and leave the "Concepts" sheet
# as the active sheet/selection (matching the authored workbook).

$wb = $excel.ActiveWorkbook

$wsMetadata = $wb.Worksheets.Item("Metadata")
$wsImported = $wb.Worksheets.Item("ImportedOntologies")
$wsConcepts = $wb.Worksheets.Item("Concepts")

# --- ImportedOntologies: point to the new emmo-inferred.ttl (beta4) URL ---
$newUrl = "https://raw.githubusercontent.com/emmo-repo/emmo-repo.github.io/master/versions/1.0.0-beta4/emmo-inferred.ttl"
$wsImported.Range("A3").Value = $newUrl
$wsImported.Hyperlinks.Add($wsImported.Range("A3"), $newUrl)

# --- Concepts: the example "Pattern" concept becomes "SpecialPattern" ---
$wsConcepts.Range("A4").Value = "SpecialPattern"

# --- Restore the per-sheet selections seen in the saved workbook ---
$wsMetadata.Activate()
$wsMetadata.Range("B20").Select()

$wsImported.Activate()
$wsImported.Range("A3").Select()

$wsConcepts.Activate()
$wsConcepts.Range("B11").Select()
